$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.164.32"
$ws.Range("E2").Value = "  +0.33%  "

# Row 3
$ws.Range("D3").Value = "1.902.82"
$ws.Range("E3").Value = "  +0.84%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.14"
$ws.Range("E5").Value = "  -0.44%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5253"
$ws.Range("E7").Value = "  +1.24%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3774"
$ws.Range("E8").Value = "  +1.35%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07255"
$ws.Range("E9").Value = "  +0.62%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.11"
$ws.Range("E10").Value = "  +0.16%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8990"
$ws.Range("E11").Value = "  -0.61%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08376"
$ws.Range("E12").Value = "  +9.30%  "

# Row 13
$ws.Range("D13").Value = "1.902.06"
$ws.Range("E13").Value = "  +0.95%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.78"
$ws.Range("E14").Value = "  -0.48%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.266"
$ws.Range("E15").Value = "  -0.03%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008618"
$ws.Range("E17").Value = "  +1.31%  "

# Row 18
$ws.Range("E18").Value = "  +1.58%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9996"

# Row 20
$ws.Range("D20").Value = "27.201.05"
$ws.Range("E20").Value = "  +0.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.061"
$ws.Range("E21").Value = "  +0.28%  "

# Row 22
$ws.Range("D22").Value = "2.146.80"
$ws.Range("E22").Value = "  +2.01%  "

# Row 23
$ws.Range("E23").Value = "  +0.36%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.432"
$ws.Range("E24").Value = "  -0.42%  "

# Row 25
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.54"
$ws.Range("E25").Value = "  +0.47%  "

# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.280"
$ws.Range("E26").Value = "  +7.06%  "

# Row 27
$ws.Range("E27").Value = "  -1.70%  "

# Row 28
$ws.Range("E28").Value = "  +0.48%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.75"
$ws.Range("E29").Value = "  +0.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.922"
$ws.Range("E30").Value = "  -0.15%  "

# Row 31
$ws.Range("E31").Value = "  -0.14%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09278"
$ws.Range("E32").Value = "  +0.79%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8129"
$ws.Range("E33").Value = "  +6.57%  "

# Row 34
$ws.Range("E34").Value = "  +0.27%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.238"
$ws.Range("E35").Value = "  +3.99%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.952"
$ws.Range("E36").Value = "  -2.18%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.350"
$ws.Range("E37").Value = "  +2.18%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.594"
$ws.Range("E38").Value = "  +1.89%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5710"
$ws.Range("E39").Value = "  +1.78%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01980"
$ws.Range("E40").Value = "  -0.71%  "

# Row 41
$ws.Range("E41").Value = "  -0.44%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.675"
$ws.Range("E42").Value = "  +1.23%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.950"
$ws.Range("E43").Value = "  +1.07%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "117.81"
$ws.Range("E44").Value = "  -0.41%  "

# Row 45
$ws.Range("E45").Value = "  +0.33%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4831"
$ws.Range("E46").Value = "  +0.83%  "

# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9999"
$ws.Range("E47").Value = "  +0.13%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.15"
$ws.Range("E48").Value = "  -0.22%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.614"
$ws.Range("E49").Value = "  +2.42%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.46"
$ws.Range("E50").Value = "  +0.91%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.64"
$ws.Range("E51").Value = "  +0.25%  "
